$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-11) before rewriting the full updated table (rows 2-16).
$ws.Range("A2:F11").ClearContents()

# Force the data range to Text format first so date-like strings in columns A/B
# (e.g. "2026-01-29") are stored as literal text, not auto-converted to date serials.
$dataRange = $ws.Range("A2:F16")
$dataRange.NumberFormat = "@"

# Write the updated data: 3 new rows inserted at the top, 1 new row inserted in the
# middle, 1 new row appended at the end; previously-existing rows shifted accordingly.
# Row 2
$ws.Cells.Item(2, 1).Value = "2026-01-29"
$ws.Cells.Item(2, 2).Value = "2026-01-28"
$ws.Cells.Item(2, 3).Value = "OpenAI"
$ws.Cells.Item(2, 4).Value = "The next chapter for AI in the EU"
$ws.Cells.Item(2, 5).Value = "EU AI의 다음 장"
$ws.Cells.Item(2, 6).Value = "https://openai.com/index/the-next-chapter-for-ai-in-the-eu"

# Row 3
$ws.Cells.Item(3, 1).Value = "2026-01-29"
$ws.Cells.Item(3, 2).Value = "2026-01-28"
$ws.Cells.Item(3, 3).Value = "OpenAI"
$ws.Cells.Item(3, 4).Value = "EMEA Youth & Wellbeing Grant"
$ws.Cells.Item(3, 5).Value = "EMEA 청소년 및 웰빙 보조금"
$ws.Cells.Item(3, 6).Value = "https://openai.com/index/emea-youth-and-wellbeing-grant"

# Row 4
$ws.Cells.Item(4, 1).Value = "2026-01-29"
$ws.Cells.Item(4, 2).Value = "2026-01-28"
$ws.Cells.Item(4, 3).Value = "OpenAI"
$ws.Cells.Item(4, 4).Value = "Keeping your data safe when an AI agent clicks a link"
$ws.Cells.Item(4, 5).Value = "AI 에이전트가 링크를 클릭할 때 데이터를 안전하게 유지"
$ws.Cells.Item(4, 6).Value = "https://openai.com/index/ai-agent-link-safety"

# Row 5
$ws.Cells.Item(5, 1).Value = "2026-01-28"
$ws.Cells.Item(5, 2).Value = "2026-01-27"
$ws.Cells.Item(5, 3).Value = "OpenAI"
$ws.Cells.Item(5, 4).Value = "PVH reimagines the future of fashion with OpenAI"
$ws.Cells.Item(5, 5).Value = "PVH는 OpenAI를 통해 패션의 미래를 재구상합니다."
$ws.Cells.Item(5, 6).Value = "https://openai.com/index/pvh-future-of-fashion"

# Row 6
$ws.Cells.Item(6, 1).Value = "2026-01-28"
$ws.Cells.Item(6, 2).Value = "2026-01-27"
$ws.Cells.Item(6, 3).Value = "OpenAI"
$ws.Cells.Item(6, 4).Value = "Powering tax donations with AI powered personalized recommendations"
$ws.Cells.Item(6, 5).Value = "AI 기반 맞춤형 추천으로 세금 기부 지원"
$ws.Cells.Item(6, 6).Value = "https://openai.com/index/trustbank"

# Row 7
$ws.Cells.Item(7, 1).Value = "2026-01-28"
$ws.Cells.Item(7, 2).Value = "2026-01-27"
$ws.Cells.Item(7, 3).Value = "OpenAI"
$ws.Cells.Item(7, 4).Value = "Introducing Prism"
$ws.Cells.Item(7, 5).Value = "프리즘 소개"
$ws.Cells.Item(7, 6).Value = "https://openai.com/index/introducing-prism"

# Row 8
$ws.Cells.Item(8, 1).Value = "2026-01-28"
$ws.Cells.Item(8, 2).Value = "2026-01-26"
$ws.Cells.Item(8, 3).Value = "OpenAI"
$ws.Cells.Item(8, 4).Value = "How Indeed uses AI to help evolve the job search"
$ws.Cells.Item(8, 5).Value = "인디드가 AI를 활용하여 구직 활동을 발전시키는 방법"
$ws.Cells.Item(8, 6).Value = "https://openai.com/index/indeed-maggie-hulce"

# Row 9
$ws.Cells.Item(9, 1).Value = "2026-01-28"
$ws.Cells.Item(9, 2).Value = "2026-01-23"
$ws.Cells.Item(9, 3).Value = "OpenAI"
$ws.Cells.Item(9, 4).Value = "Unrolling the Codex agent loop"
$ws.Cells.Item(9, 5).Value = "Codex 에이전트 루프 풀기"
$ws.Cells.Item(9, 6).Value = "https://openai.com/index/unrolling-the-codex-agent-loop"

# Row 10
$ws.Cells.Item(10, 1).Value = "2026-01-28"
$ws.Cells.Item(10, 2).Value = "2026-01-22"
$ws.Cells.Item(10, 3).Value = "OpenAI"
$ws.Cells.Item(10, 4).Value = "Scaling PostgreSQL to power 800 million ChatGPT users"
$ws.Cells.Item(10, 5).Value = "8억 명의 ChatGPT 사용자를 지원하기 위해 PostgreSQL 확장"
$ws.Cells.Item(10, 6).Value = "https://openai.com/index/scaling-postgresql"

# Row 11
$ws.Cells.Item(11, 1).Value = "2026-01-28"
$ws.Cells.Item(11, 2).Value = "2026-01-22"
$ws.Cells.Item(11, 3).Value = "OpenAI"
$ws.Cells.Item(11, 4).Value = "Inside Praktika's conversational approach to language learning"
$ws.Cells.Item(11, 5).Value = "Praktika의 언어 학습에 대한 대화식 접근 방식 살펴보기"
$ws.Cells.Item(11, 6).Value = "https://openai.com/index/praktika"

# Row 12
$ws.Cells.Item(12, 1).Value = "2026-01-28"
$ws.Cells.Item(12, 2).Value = "2026-01-22"
$ws.Cells.Item(12, 3).Value = "OpenAI"
$ws.Cells.Item(12, 4).Value = "Inside GPT-5 for Work: How Businesses Use GPT-5"
$ws.Cells.Item(12, 5).Value = "업무용 GPT-5 내부: 기업이 GPT-5를 사용하는 방법"
$ws.Cells.Item(12, 6).Value = "https://openai.com/business/guides-and-resources/chatgpt-usage-and-adoption-patterns-at-work"

# Row 13
$ws.Cells.Item(13, 1).Value = "2026-01-29"
$ws.Cells.Item(13, 2).Value = "2026-01-21"
$ws.Cells.Item(13, 3).Value = "OpenAI"
$ws.Cells.Item(13, 4).Value = "How countries can end the capability overhang"
$ws.Cells.Item(13, 5).Value = "국가가 역량 과잉을 어떻게 끝낼 수 있는가"
$ws.Cells.Item(13, 6).Value = "https://openai.com/index/how-countries-can-end-the-capability-overhang"

# Row 14
$ws.Cells.Item(14, 1).Value = "2026-01-28"
$ws.Cells.Item(14, 2).Value = "2026-01-21"
$ws.Cells.Item(14, 3).Value = "OpenAI"
$ws.Cells.Item(14, 4).Value = "How Higgsfield turns simple ideas into cinematic social videos"
$ws.Cells.Item(14, 5).Value = "Higgsfield가 단순한 아이디어를 영화 같은 소셜 비디오로 바꾸는 방법"
$ws.Cells.Item(14, 6).Value = "https://openai.com/index/higgsfield"

# Row 15
$ws.Cells.Item(15, 1).Value = "2026-01-28"
$ws.Cells.Item(15, 2).Value = "2026-01-21"
$ws.Cells.Item(15, 3).Value = "OpenAI"
$ws.Cells.Item(15, 4).Value = "Introducing Edu for Countries"
$ws.Cells.Item(15, 5).Value = "국가별 교육 소개"
$ws.Cells.Item(15, 6).Value = "https://openai.com/index/edu-for-countries"

# Row 16
$ws.Cells.Item(16, 1).Value = "2026-01-29"
$ws.Cells.Item(16, 2).Value = "2026-01-20"
$ws.Cells.Item(16, 3).Value = "OpenAI"
$ws.Cells.Item(16, 4).Value = "Horizon 1000: Advancing AI for primary healthcare"
$ws.Cells.Item(16, 5).Value = "Horizon 1000: 1차 의료를 위한 AI 발전"
$ws.Cells.Item(16, 6).Value = "https://openai.com/index/horizon-1000"

# Reset the style back to Normal (default, no explicit style index) to match the
# plain formatting of the rest of the data rows, while keeping values as text.
$dataRange.Style = "Normal"

Write-Host "OpenAI news table updated (rows 2-16)."
